$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'45.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.582"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.83%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08341"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.79%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.117"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.35%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9774"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.59%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.69%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'5.64%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'1.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-3.45%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.70%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04668"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.36%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.59%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001288"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.38%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005882"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.63%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.376"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.29%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.450"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.28%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-2.55%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.44%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2784"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'11.30%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04173"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.31%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001294"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.75%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004592"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'5.13%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'8.69%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.21%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.80%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05740"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.82%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007888"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.14%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1431"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.86%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007518"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.91%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002101"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.36%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008514"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.56%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3369"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007128"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.17%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.33%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'0.39%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003528"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.64%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.59%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.33%"
$ws.Range("E51").Style = "Normal"
